$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the test-data row (row 5) with the new POC values used against Stage A.
$ws.Range("A5").Value = "AH_healthsys_poc1"
$ws.Range("B5").Value = "AH_healthsys_poc1"

$ws.Range("C5").Value = "Facility_POC_2"
$ws.Range("D5").Value = 4112019

$ws.Range("E5").Value = "Pharmacy_POC_2"
$ws.Range("F5").Value = "Pharmacy_POC_2"

# Reflect the editor's final cursor position at save time.
[void]$ws.Range("G18").Select()
